$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H so the old "Total Year Up" (H) shifts to I,
# making room for the new "To Year High Profit" column at H.
$ws.Columns.Item(8).Insert()

# New header for inserted column H; old H1 ("Total Year Up") already shifted to I1.
$ws.Cells.Item(1, 8).Value = "To Year High Profit"

# Rewrite all data rows (values updated + reordered by Total Year Down desc; one new row added).
# Row 2: ADANIGREEN
$ws.Cells.Item(2, 1).Value = "INE364U01010"
$ws.Cells.Item(2, 2).Value = "ADANIGREEN"
$ws.Cells.Item(2, 3).Value = "Power Generation"
$ws.Cells.Item(2, 4).Value = 1039
$ws.Cells.Item(2, 5).Value = 2174.1
$ws.Cells.Item(2, 6).Value = 870.25
$ws.Cells.Item(2, 7).Value = 52.21010993054598
$ws.Cells.Item(2, 8).Value = 109.2492781520693
$ws.Cells.Item(2, 9).Value = 19.3909796035622

# Row 3: ADANIPOWER
$ws.Cells.Item(3, 1).Value = "INE814H01011"
$ws.Cells.Item(3, 2).Value = "ADANIPOWER"
$ws.Cells.Item(3, 3).Value = "Integrated Power Utilities"
$ws.Cells.Item(3, 4).Value = 519
$ws.Cells.Item(3, 5).Value = 895.85
$ws.Cells.Item(3, 6).Value = 432
$ws.Cells.Item(3, 7).Value = 42.06619411731874
$ws.Cells.Item(3, 8).Value = 72.61078998073218
$ws.Cells.Item(3, 9).Value = 20.13888888888889

# Row 4: INDUSINDBK
$ws.Cells.Item(4, 1).Value = "INE095A01012"
$ws.Cells.Item(4, 2).Value = "INDUSINDBK"
$ws.Cells.Item(4, 3).Value = "Private Sector Bank"
$ws.Cells.Item(4, 4).Value = 996
$ws.Cells.Item(4, 5).Value = 1694.5
$ws.Cells.Item(4, 6).Value = 926.45
$ws.Cells.Item(4, 7).Value = 41.22159929182649
$ws.Cells.Item(4, 8).Value = 70.13052208835342
$ws.Cells.Item(4, 9).Value = 7.507150952560848

# Row 5: ADANIENSOL
$ws.Cells.Item(5, 1).Value = "INE931S01010"
$ws.Cells.Item(5, 2).Value = "ADANIENSOL"
$ws.Cells.Item(5, 3).Value = "Power Distribution"
$ws.Cells.Item(5, 4).Value = 821
$ws.Cells.Item(5, 5).Value = 1348
$ws.Cells.Item(5, 6).Value = 588
$ws.Cells.Item(5, 7).Value = 39.09495548961425
$ws.Cells.Item(5, 8).Value = 64.19001218026796
$ws.Cells.Item(5, 9).Value = 39.62585034013606

# Row 6: ATGL
$ws.Cells.Item(6, 1).Value = "INE399L01023"
$ws.Cells.Item(6, 2).Value = "ATGL"
$ws.Cells.Item(6, 3).Value = "LPG/CNG/PNG/LNG Supplier"
$ws.Cells.Item(6, 4).Value = 725.8
$ws.Cells.Item(6, 5).Value = 1190
$ws.Cells.Item(6, 6).Value = 545.75
$ws.Cells.Item(6, 7).Value = 39.00840336134455
$ws.Cells.Item(6, 8).Value = 63.95701295122625
$ws.Cells.Item(6, 9).Value = 32.99129638112688

# Row 7: TATAMOTORS
$ws.Cells.Item(7, 1).Value = "INE155A01022"
$ws.Cells.Item(7, 2).Value = "TATAMOTORS"
$ws.Cells.Item(7, 3).Value = "Passenger Cars & Utility Vehicles"
$ws.Cells.Item(7, 4).Value = 789
$ws.Cells.Item(7, 5).Value = 1179
$ws.Cells.Item(7, 6).Value = 717.7
$ws.Cells.Item(7, 7).Value = 33.07888040712468
$ws.Cells.Item(7, 8).Value = 49.42965779467681
$ws.Cells.Item(7, 9).Value = 9.934513027727455

# Row 8: IRFC
$ws.Cells.Item(8, 1).Value = "INE053F01010"
$ws.Cells.Item(8, 2).Value = "IRFC"
$ws.Cells.Item(8, 3).Value = "Financial Institution"
$ws.Cells.Item(8, 4).Value = 153.56
$ws.Cells.Item(8, 5).Value = 229
$ws.Cells.Item(8, 6).Value = 99
$ws.Cells.Item(8, 7).Value = 32.94323144104803
$ws.Cells.Item(8, 8).Value = 49.12737692107319
$ws.Cells.Item(8, 9).Value = 55.11111111111111

# Row 9: HEROMOTOCO
$ws.Cells.Item(9, 1).Value = "INE158A01026"
$ws.Cells.Item(9, 2).Value = "HEROMOTOCO"
$ws.Cells.Item(9, 3).Value = "2/3 Wheelers"
$ws.Cells.Item(9, 4).Value = 4240.05
$ws.Cells.Item(9, 5).Value = 6246.25
$ws.Cells.Item(9, 6).Value = 3929.85
$ws.Cells.Item(9, 7).Value = 32.11847108264959
$ws.Cells.Item(9, 8).Value = 47.31547977028572
$ws.Cells.Item(9, 9).Value = 7.893431046986543

# Row 10: ADANIENT
$ws.Cells.Item(10, 1).Value = "INE423A01024"
$ws.Cells.Item(10, 2).Value = "ADANIENT"
$ws.Cells.Item(10, 3).Value = "Trading - Minerals"
$ws.Cells.Item(10, 4).Value = 2562
$ws.Cells.Item(10, 5).Value = 3743.9
$ws.Cells.Item(10, 6).Value = 2025
$ws.Cells.Item(10, 7).Value = 31.56868506103262
$ws.Cells.Item(10, 8).Value = 46.13192818110852
$ws.Cells.Item(10, 9).Value = 26.51851851851852

# Row 11: BHEL
$ws.Cells.Item(11, 1).Value = "INE257A01026"
$ws.Cells.Item(11, 2).Value = "BHEL"
$ws.Cells.Item(11, 3).Value = "Heavy Electrical Equipment"
$ws.Cells.Item(11, 4).Value = 229.99
$ws.Cells.Item(11, 5).Value = 335.35
$ws.Cells.Item(11, 6).Value = 191.85
$ws.Cells.Item(11, 7).Value = 31.41792157447443
$ws.Cells.Item(11, 8).Value = 45.81068742119223
$ws.Cells.Item(11, 9).Value = 19.88011467292157

# Row 12: ASIANPAINT
$ws.Cells.Item(12, 1).Value = "INE021A01026"
$ws.Cells.Item(12, 2).Value = "ASIANPAINT"
$ws.Cells.Item(12, 3).Value = "Paints"
$ws.Cells.Item(12, 4).Value = 2333.8
$ws.Cells.Item(12, 5).Value = 3401.25
$ws.Cells.Item(12, 6).Value = 2256.5
$ws.Cells.Item(12, 7).Value = 31.3840499816244
$ws.Cells.Item(12, 8).Value = 45.73870940097693
$ws.Cells.Item(12, 9).Value = 3.425659206736098

# Row 13: IRCTC
$ws.Cells.Item(13, 1).Value = "INE335Y01020"
$ws.Cells.Item(13, 2).Value = "IRCTC"
$ws.Cells.Item(13, 3).Value = "Tour Travel Related Services"
$ws.Cells.Item(13, 4).Value = 795.2
$ws.Cells.Item(13, 5).Value = 1138.9
$ws.Cells.Item(13, 6).Value = 765.1
$ws.Cells.Item(13, 7).Value = 30.17824216349109
$ws.Cells.Item(13, 8).Value = 43.2218309859155
$ws.Cells.Item(13, 9).Value = 3.934126258005488

# Row 14: NHPC
$ws.Cells.Item(14, 1).Value = "INE848E01016"
$ws.Cells.Item(14, 2).Value = "NHPC"
$ws.Cells.Item(14, 3).Value = "Power Generation"
$ws.Cells.Item(14, 4).Value = 82.73
$ws.Cells.Item(14, 5).Value = 118.4
$ws.Cells.Item(14, 6).Value = 67
$ws.Cells.Item(14, 7).Value = 30.12668918918919
$ws.Cells.Item(14, 8).Value = 43.11616100568114
$ws.Cells.Item(14, 9).Value = 23.47761194029852

# Row 15: IOC
$ws.Cells.Item(15, 1).Value = "INE242A01010"
$ws.Cells.Item(15, 2).Value = "IOC"
$ws.Cells.Item(15, 3).Value = "Refineries & Marketing"
$ws.Cells.Item(15, 4).Value = 138.05
$ws.Cells.Item(15, 5).Value = 196.8
$ws.Cells.Item(15, 6).Value = 128.5
$ws.Cells.Item(15, 7).Value = 29.85264227642276
$ws.Cells.Item(15, 8).Value = 42.55704454907643
$ws.Cells.Item(15, 9).Value = 7.43190661478601

# Row 16: BAJAJ-AUTO
$ws.Cells.Item(16, 1).Value = "INE917I01010"
$ws.Cells.Item(16, 2).Value = "BAJAJ-AUTO"
$ws.Cells.Item(16, 3).Value = "2/3 Wheelers"
$ws.Cells.Item(16, 4).Value = 8965
$ws.Cells.Item(16, 5).Value = 12774
$ws.Cells.Item(16, 6).Value = 6604
$ws.Cells.Item(16, 7).Value = 29.81838108658213
$ws.Cells.Item(16, 8).Value = 42.48745119910764
$ws.Cells.Item(16, 9).Value = 35.75105996365839

# Row 17: UNIONBANK
$ws.Cells.Item(17, 1).Value = "INE692A01016"
$ws.Cells.Item(17, 2).Value = "UNIONBANK"
$ws.Cells.Item(17, 3).Value = "Public Sector Bank"
$ws.Cells.Item(17, 4).Value = 124
$ws.Cells.Item(17, 5).Value = 172.5
$ws.Cells.Item(17, 6).Value = 106.68
$ws.Cells.Item(17, 7).Value = 28.11594202898551
$ws.Cells.Item(17, 8).Value = 39.11290322580645
$ws.Cells.Item(17, 9).Value = 16.23547056617922

# Row 18: GODREJCP
$ws.Cells.Item(18, 1).Value = "INE102D01028"
$ws.Cells.Item(18, 2).Value = "GODREJCP"
$ws.Cells.Item(18, 3).Value = "Personal Care"
$ws.Cells.Item(18, 4).Value = 1115
$ws.Cells.Item(18, 5).Value = 1541.85
$ws.Cells.Item(18, 6).Value = 1055.05
$ws.Cells.Item(18, 7).Value = 27.68427538346791
$ws.Cells.Item(18, 8).Value = 38.28251121076232
$ws.Cells.Item(18, 9).Value = 5.682195156627645

# Row 19: COALINDIA
$ws.Cells.Item(19, 1).Value = "INE522F01014"
$ws.Cells.Item(19, 2).Value = "COALINDIA"
$ws.Cells.Item(19, 3).Value = "Coal"
$ws.Cells.Item(19, 4).Value = 394.1
$ws.Cells.Item(19, 5).Value = 543.55
$ws.Cells.Item(19, 6).Value = 368
$ws.Cells.Item(19, 7).Value = 27.49517063747584
$ws.Cells.Item(19, 8).Value = 37.92184724689163
$ws.Cells.Item(19, 9).Value = 7.092391304347823

# Row 20: DMART
$ws.Cells.Item(20, 1).Value = "INE192R01011"
$ws.Cells.Item(20, 2).Value = "DMART"
$ws.Cells.Item(20, 3).Value = "Diversified Retail"
$ws.Cells.Item(20, 4).Value = 4011.9
$ws.Cells.Item(20, 5).Value = 5484.85
$ws.Cells.Item(20, 6).Value = 3399
$ws.Cells.Item(20, 7).Value = 26.85488208428673
$ws.Cells.Item(20, 8).Value = 36.71452429023656
$ws.Cells.Item(20, 9).Value = 18.03177405119154

# Row 21: ZYDUSLIFE
$ws.Cells.Item(21, 1).Value = "INE010B01027"
$ws.Cells.Item(21, 2).Value = "ZYDUSLIFE"
$ws.Cells.Item(21, 3).Value = "Pharmaceuticals"
$ws.Cells.Item(21, 4).Value = 971.75
$ws.Cells.Item(21, 5).Value = 1324.3
$ws.Cells.Item(21, 6).Value = 686
$ws.Cells.Item(21, 7).Value = 26.62161141735256
$ws.Cells.Item(21, 8).Value = 36.27990738358631
$ws.Cells.Item(21, 9).Value = 41.65451895043732

# Row 22: MOTHERSON
$ws.Cells.Item(22, 1).Value = "INE775A01035"
$ws.Cells.Item(22, 2).Value = "MOTHERSON"
$ws.Cells.Item(22, 3).Value = "Auto Components & Equipments"
$ws.Cells.Item(22, 4).Value = 159.45
$ws.Cells.Item(22, 5).Value = 216.99
$ws.Cells.Item(22, 6).Value = 101.35
$ws.Cells.Item(22, 7).Value = 26.51735103000139
$ws.Cells.Item(22, 8).Value = 36.08654750705551
$ws.Cells.Item(22, 9).Value = 57.32609768130241

# Row 23: ADANIPORTS
$ws.Cells.Item(23, 1).Value = "INE742F01042"
$ws.Cells.Item(23, 2).Value = "ADANIPORTS"
$ws.Cells.Item(23, 3).Value = "Port & Port services"
$ws.Cells.Item(23, 4).Value = 1197.8
$ws.Cells.Item(23, 5).Value = 1621.4
$ws.Cells.Item(23, 6).Value = 995.65
$ws.Cells.Item(23, 7).Value = 26.12557049463427
$ws.Cells.Item(23, 8).Value = 35.36483553180834
$ws.Cells.Item(23, 9).Value = 20.3033194395621

# Row 24: HAL
$ws.Cells.Item(24, 1).Value = "INE066F01020"
$ws.Cells.Item(24, 2).Value = "HAL"
$ws.Cells.Item(24, 3).Value = "Aerospace & Defense"
$ws.Cells.Item(24, 4).Value = 4203
$ws.Cells.Item(24, 5).Value = 5674.75
$ws.Cells.Item(24, 6).Value = 2763
$ws.Cells.Item(24, 7).Value = 25.93506321864399
$ws.Cells.Item(24, 8).Value = 35.0166547704021
$ws.Cells.Item(24, 9).Value = 52.11726384364821

# Row 25: ABB
$ws.Cells.Item(25, 1).Value = "INE117A01022"
$ws.Cells.Item(25, 2).Value = "ABB"
$ws.Cells.Item(25, 3).Value = "Heavy Electrical Equipment"
$ws.Cells.Item(25, 4).Value = 6782
$ws.Cells.Item(25, 5).Value = 9149.950000000001
$ws.Cells.Item(25, 6).Value = 4340.3
$ws.Cells.Item(25, 7).Value = 25.87937639003492
$ws.Cells.Item(25, 8).Value = 34.91521675022118
$ws.Cells.Item(25, 9).Value = 56.25647996682257

# Row 26: LICI
$ws.Cells.Item(26, 1).Value = "INE0J1Y01017"
$ws.Cells.Item(26, 2).Value = "LICI"
$ws.Cells.Item(26, 3).Value = "Life Insurance"
$ws.Cells.Item(26, 4).Value = 908
$ws.Cells.Item(26, 5).Value = 1222
$ws.Cells.Item(26, 6).Value = 819.3
$ws.Cells.Item(26, 7).Value = 25.69558101472995
$ws.Cells.Item(26, 8).Value = 34.58149779735682
$ws.Cells.Item(26, 9).Value = 10.82631514707677

# Row 27: TATACONSUM
$ws.Cells.Item(27, 1).Value = "INE192A01025"
$ws.Cells.Item(27, 2).Value = "TATACONSUM"
$ws.Cells.Item(27, 3).Value = "Tea & Coffee"
$ws.Cells.Item(27, 4).Value = 936
$ws.Cells.Item(27, 5).Value = 1256.44
$ws.Cells.Item(27, 6).Value = 882.9
$ws.Cells.Item(27, 7).Value = 25.50380439973258
$ws.Cells.Item(27, 8).Value = 34.23504273504274
$ws.Cells.Item(27, 9).Value = 6.014271151885842

# Row 28: PNB
$ws.Cells.Item(28, 1).Value = "INE160A01022"
$ws.Cells.Item(28, 2).Value = "PNB"
$ws.Cells.Item(28, 3).Value = "Public Sector Bank"
$ws.Cells.Item(28, 4).Value = 106.53
$ws.Cells.Item(28, 5).Value = 142.9
$ws.Cells.Item(28, 6).Value = 92.40000000000001
$ws.Cells.Item(28, 7).Value = 25.45136459062282
$ws.Cells.Item(28, 8).Value = 34.14061766638505
$ws.Cells.Item(28, 9).Value = 15.29220779220779

# Row 29: BRITANNIA
$ws.Cells.Item(29, 1).Value = "INE216A01030"
$ws.Cells.Item(29, 2).Value = "BRITANNIA"
$ws.Cells.Item(29, 3).Value = "Packaged Foods"
$ws.Cells.Item(29, 4).Value = 4839.85
$ws.Cells.Item(29, 5).Value = 6469.9
$ws.Cells.Item(29, 6).Value = 4641
$ws.Cells.Item(29, 7).Value = 25.19436158209554
$ws.Cells.Item(29, 8).Value = 33.67976280256617
$ws.Cells.Item(29, 9).Value = 4.284636931695762

# Row 30: SBILIFE
$ws.Cells.Item(30, 1).Value = "INE123W01016"
$ws.Cells.Item(30, 2).Value = "SBILIFE"
$ws.Cells.Item(30, 3).Value = "Life Insurance"
$ws.Cells.Item(30, 4).Value = 1450
$ws.Cells.Item(30, 5).Value = 1936
$ws.Cells.Item(30, 6).Value = 1307.7
$ws.Cells.Item(30, 7).Value = 25.10330578512396
$ws.Cells.Item(30, 8).Value = 33.51724137931033
$ws.Cells.Item(30, 9).Value = 10.88170069587826

# Row 31: TATASTEEL
$ws.Cells.Item(31, 1).Value = "INE081A01020"
$ws.Cells.Item(31, 2).Value = "TATASTEEL"
$ws.Cells.Item(31, 3).Value = "Iron & Steel"
$ws.Cells.Item(31, 4).Value = 138.33
$ws.Cells.Item(31, 5).Value = 184.6
$ws.Cells.Item(31, 6).Value = 128.2
$ws.Cells.Item(31, 7).Value = 25.06500541711808
$ws.Cells.Item(31, 8).Value = 33.44899877105472
$ws.Cells.Item(31, 9).Value = 7.901716068642761

# Row 32: ONGC
$ws.Cells.Item(32, 1).Value = "INE213A01029"
$ws.Cells.Item(32, 2).Value = "ONGC"
$ws.Cells.Item(32, 3).Value = "Oil Exploration & Production"
$ws.Cells.Item(32, 4).Value = 258.65
$ws.Cells.Item(32, 5).Value = 345
$ws.Cells.Item(32, 6).Value = 204.5
$ws.Cells.Item(32, 7).Value = 25.02898550724638
$ws.Cells.Item(32, 8).Value = 33.38488304658807
$ws.Cells.Item(32, 9).Value = 26.47921760391196

Write-Output "done"